$d = $word.ActiveDocument

$pkgPrefix = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:body>'
$pkgSuffix = '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# --- Edit 1: fill empty paragraph 13 (before "funktioner vi valgt...") ---
$body1 = '<w:p><w:pPr><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve">Først kalder vi vores class og giver den et navn, derefter laver vi en list som skal indeholde en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>string</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>(tekst),</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:br/><w:t xml:space="preserve">så laver vi en </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>int</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>(tal), kalder vores class og giver den nogen paramenter,</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:br/><w:t xml:space="preserve">give vores kode til at logge ind og til at komme ind i </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>enable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve">, så lave vi nogen </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>commands</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>,</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:br/><w:t xml:space="preserve">til sidst lave vi et loop som viser alle de linjer som er vi vores startup </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>config</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve">.  </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr></w:pPr></w:p>'
$p1 = $d.Paragraphs.Item(13)
$r1 = $p1.Range
$r1.InsertXML($pkgPrefix + $body1 + $pkgSuffix)

# --- Edit 2: append to paragraph 18 (try/catch, was 17, shifted by +1 after edit 1) ---
$body2 = '<w:p><w:pPr><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve">På billedet her kan vi se de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>try</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>catch</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve"> vi har brugt.</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:br/><w:t xml:space="preserve">Tjekker om det du skrive er i </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>den rigtige format</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>.</w:t></w:r></w:p>'
$p2 = $d.Paragraphs.Item(18)
$r2 = $p2.Range
$r2.InsertXML($pkgPrefix + $body2 + $pkgSuffix)

# --- Edit 3: replace paragraphs 19+20 (image, if-statements; was 18+19, shifted by +1) ---
$body3 = '<w:p><w:pPr><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="da-DK"/></w:rPr><w:lastRenderedPageBreak/><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="38D914AB" wp14:editId="019C8576"><wp:extent cx="3765744" cy="4375375"/><wp:effectExtent l="0" t="0" r="6350" b="6350"/><wp:docPr id="5" name="Billede 5"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="5" name="if elseif statment.PNG"/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId10"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="3765744" cy="4375375"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>På bille</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve">det her er vores </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>if</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve"> statements.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>Her ser du hvad der sker ved de forskellige funktioner. F.eks. hvis du t</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>r</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>ykker 1</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>(som er lig med r1)</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve"> vil den gå i gang med vores Router konfiguration, og der efter vil der komme en linje med en tekst hvor du kan tykke en knap af eget valg for at exit router funktionen eller hvis du trykker 2</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t>(som lig med sw1)</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve"> så vil den gå videre med vores switch konfiguration hvor du så efter den er færdig igen kan trykke på en kap efter eget valg for at exit.</w:t></w:r><w:r><w:rPr><w:lang w:eastAsia="da-DK"/></w:rPr><w:t xml:space="preserve"> Hvis du så trykker 3(som er lig med begge funktioner) vil den først gå i gang med vores router funktion, derefter skriver den så ”så går vi videre med SW1” hvor du så skal trykke på en knap efter eget valg for at gå videre med vores switch konfiguration, hvor du så igen skal trykke på en knap efter eget valg for at exit.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'
$pA = $d.Paragraphs.Item(19)
$pB = $d.Paragraphs.Item(20)
$r3 = $d.Range($pA.Range.Start, $pB.Range.End)
$r3.InsertXML($pkgPrefix + $body3 + $pkgSuffix)

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
